# Fix mojibake: "Â±" (U+00C2 U+00B1, a double UTF-8-encoded "±") -> "±" (U+00B1)
# Affects the data cells (columns B:H, rows 2:17) on the f1_score, training_time
# and test_time sheets.

$wb = $excel.ActiveWorkbook

$bad  = [string]([char]0x00C2) + [string]([char]0x00B1)
$good = [string]([char]0x00B1)

$sheetNames = @("f1_score", "training_time", "test_time")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 17; $row++) {
        for ($col = 2; $col -le 8; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $text = $cell.Text
            if ($text.Contains($bad)) {
                $cell.Value = $text.Replace($bad, $good)
            }
        }
    }
}
